$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 21 formatting: turn it into a section-divider row like row 17 ---
# Copy the "divider" formats (border-bottom, no alignment override) from row 17
# onto row 21's cells, without touching their values.
$ws.Cells.Item(17, 2).Copy() | Out-Null
$ws.Cells.Item(21, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> A21 (was empty)
$ws.Cells.Item(21, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> B21 (s=4 -> s=6)

$ws.Cells.Item(17, 3).Copy() | Out-Null
$ws.Cells.Item(21, 3).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> C21 (s=5 -> s=7)
$ws.Cells.Item(21, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> D21 (s=5 -> s=7)
$ws.Cells.Item(21, 5).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> E21 (s=5 -> s=7)

$excel.CutCopyMode = 0

# --- Add new rows 22 & 23 with the same layout as the other data rows ---
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 43.2

# Row 22: full data row (filename, line number, english, translated, converted)
# Write the string cells in the same order the new shared strings appear in the
# diff: English text, filename, translated text, converted text, then the
# row-23 filename.
$ws.Cells.Item(22, 3).Value2 = " Oooh, yeah... Soaking in the Hot\nSpring is so relaxing…"
$ws.Cells.Item(22, 1).Value2 = "SCRIPT/P02P01A/us0102.ssb"
$ws.Cells.Item(22, 4).Value2 = " Ооо, да... Как хорошо купаться в\nГорячих Источниках..."
$ws.Cells.Item(22, 5).Value2 = " Ïïï, äà... Ëàë öïñïšï ëôðàóûòÿ â\nÃïñÿœéö Éòóïœîéëàö…"
$ws.Cells.Item(23, 1).Value2 = "SCRIPT/P02P01A/us3102.ssb"

$ws.Cells.Item(22, 2).Value2 = 18

# Copy formats for the new rows from row 20 (a normal data row) so the cell
# styles (s=4 for A/B, s=5 for C/D/E) match.
$ws.Cells.Item(20, 1).Copy() | Out-Null
$ws.Cells.Item(22, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(23, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(20, 2).Copy() | Out-Null
$ws.Cells.Item(22, 2).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(20, 3).Copy() | Out-Null
$ws.Cells.Item(22, 3).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(20, 4).Copy() | Out-Null
$ws.Cells.Item(22, 4).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(20, 5).Copy() | Out-Null
$ws.Cells.Item(22, 5).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Update the view state to match where the user scrolled/selected ---
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("E22").Select() | Out-Null
